# update data on Feb-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Five existing (previously discharged-pending) cases now show a
# "Discharged" status (column K) with a DischargeDate of Feb-19
# (column N): rows 2, 33, 34, 45, 56.
# ---------------------------------------------------------------
$dischargedRows = @(2, 33, 34, 45, 56)
foreach ($r in $dischargedRows) {
    $ws.Cells.Item($r, 11).Value = "Discharged"   # column K - Status
    $ws.Cells.Item($r, 14).Value = "Feb-19"       # column N - DischargeDate
}

# ---------------------------------------------------------------
# Row 81 (Case 80) - corrected coordinates / extra info
# ---------------------------------------------------------------
$ws.Cells.Item(81, 2).Value = 1.3965160000000001     # B - lat
$ws.Cells.Item(81, 3).Value = 103.87891999999999     # C - lon
$ws.Cells.Item(81, 8).Value = "Fernvale Link"        # H - Stay
$ws.Cells.Item(81, 9).Value = "NUH, GP Clinic"       # I - Visited
$ws.Cells.Item(81, 11).Value = "NUH"                 # K - Status
$ws.Cells.Item(81, 13).Value = "Feb-04"              # M - SymtomDate
$ws.Cells.Item(81, 16).Value = "NUH administration"  # P - Occupation

# ---------------------------------------------------------------
# Row 82 (Case 81) - corrected coordinates / extra info
# ---------------------------------------------------------------
$ws.Cells.Item(82, 2).Value = 1.280619               # B - lat
$ws.Cells.Item(82, 3).Value = 103.82365299999999     # C - lon
$ws.Cells.Item(82, 8).Value = "Lower Delta Road"     # H - Stay
$ws.Cells.Item(82, 13).Value = "Feb-16"              # M - SymtomDate

# ---------------------------------------------------------------
# New row 83 (Case 82)
# ---------------------------------------------------------------
$ws.Cells.Item(83, 1).Value = 82                      # A - Case
$ws.Cells.Item(83, 2).Value = 1.348231                # B - lat
$ws.Cells.Item(83, 3).Value = 103.72327               # C - lon
$ws.Cells.Item(83, 4).Value = "Feb-18"                # D - ConfirmDate
$ws.Cells.Item(83, 5).Value = 57                      # E - Age
$ws.Cells.Item(83, 6).Value = "Female"                # F - Gender
$ws.Cells.Item(83, 7).Value = "Singapore"             # G - WhereInfected
$ws.Cells.Item(83, 8).Value = "Jurong West Street 41" # H - Stay
$ws.Cells.Item(83, 9).Value = "Jurong Polyclinic, GP Clinic, NTFGH emergency department"  # I - Visited
$ws.Cells.Item(83, 11).Value = "NTFGH"                # K - Status
$ws.Cells.Item(83, 12).Value = "Singaporean"          # L - Nationality
$ws.Cells.Item(83, 13).Value = "Feb-09"               # M - SymtomDate
foreach ($col in @(1, 2, 5, 6, 7, 9, 12)) {
    $ws.Cells.Item(83, $col).Font.Size = 16
}

# ---------------------------------------------------------------
# New row 84 (Case 83)
# ---------------------------------------------------------------
$ws.Cells.Item(84, 1).Value = 83                      # A - Case
$ws.Cells.Item(84, 2).Value = 1.322101                # B - lat
$ws.Cells.Item(84, 3).Value = 103.847257              # C - lon
$ws.Cells.Item(84, 4).Value = "Feb-19"                # D - ConfirmDate
$ws.Cells.Item(84, 5).Value = 54                      # E - Age
$ws.Cells.Item(84, 6).Value = "Male"                  # F - Gender
$ws.Cells.Item(84, 7).Value = "Singapore"             # G - WhereInfected
$ws.Cells.Item(84, 9).Value = "Malaysia"              # I - Visited
$ws.Cells.Item(84, 11).Value = "NCID"                 # K - Status
$ws.Cells.Item(84, 12).Value = "Singaporean"          # L - Nationality
$ws.Cells.Item(84, 15).Value = "The Life Church and Missions Singapore"  # O - Cluster
foreach ($col in @(1, 2, 5, 6, 7, 9, 11, 12)) {
    $ws.Cells.Item(84, $col).Font.Size = 16
}

# ---------------------------------------------------------------
# New row 85 (Case 84)
# ---------------------------------------------------------------
$ws.Cells.Item(85, 1).Value = 84                      # A - Case
$ws.Cells.Item(85, 2).Value = 1.3221160000000001      # B - lat
$ws.Cells.Item(85, 3).Value = 103.847244              # C - lon
$ws.Cells.Item(85, 4).Value = "Feb-19"                # D - ConfirmDate
$ws.Cells.Item(85, 5).Value = 35                      # E - Age
$ws.Cells.Item(85, 6).Value = "Female"                # F - Gender
$ws.Cells.Item(85, 7).Value = "Singapore"             # G - WhereInfected
$ws.Cells.Item(85, 10).Value = "66"                   # J - LinkedTo
$ws.Cells.Item(85, 11).Value = "NCID"                 # K - Status
$ws.Cells.Item(85, 12).Value = "Singaporean"          # L - Nationality
$ws.Cells.Item(85, 15).Value = "Grace Assembly of God"  # O - Cluster
foreach ($col in @(1, 2, 5, 6, 7, 11, 12)) {
    $ws.Cells.Item(85, $col).Font.Size = 16
}

# ---------------------------------------------------------------
# Row heights for the three newly-added rows, matching the rest of
# the table (ht="21").
# ---------------------------------------------------------------
$ws.Rows.Item(83).RowHeight = 21
$ws.Rows.Item(84).RowHeight = 21
$ws.Rows.Item(85).RowHeight = 21

# ---------------------------------------------------------------
# Update the visible selection to match the new bottom of the table.
# ---------------------------------------------------------------
[void]$ws.Range("C85").Select()
